# PowerShell COM-interop script to update PtX_demand_LV Outputs sheet
# Adds "Fossil Gases" and "Fossil Liquids" rows to each year block (2030/2040/2050)
# and (re)populates values, including the newly-populated "Pass Aviation" column,
# to match the corrected category breakdown described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the 6 new rows (2 per year block), working bottom-to-top so row numbers
# used for each Insert() call stay valid as later inserts shift things below them.
$insertRows = @(29, 27, 19, 17, 9, 7)
foreach ($r in $insertRows) {
    $ws.Rows.Item($r).Insert()
}

# Full target data for rows 2-37, columns A-K (FuelGroup, Year, Iron & steel,
# Chemicals, Non-metallic minerals, Pass Road, Pass Rail, Pass Aviation,
# Freight Road, Freight Rail, Maritime).
$data = @(
    @("Hydrogen", 2030, $null, $null, $null, 0.00010994556158042990, $null, 0.00000000029940609626369042, 0.000038603887076174082, $null, $null),
    @("Methanol", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @("Ammonia", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @("Synthetic Gases", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @("Biogenic Gases", 2030, $null, $null, 0.00012926908324502090, 0.000033056771124520503, $null, $null, 0.0000084332502379877768, $null, $null),
    @("Fossil Gases", 2030, $null, $null, $null, 0.00042332070346830431, $null, $null, 0.000031032431091703537, $null, $null),
    @("Synthetic Liquids", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @("Biogenic Liquids", 2030, $null, $null, $null, 0.0019815148069229219, 0.000011536033730079161, 0.00081920397375579996, 0.0014481522324509999, 0.00016413634951250000, 0.0038832898971645690),
    @("Fossil Liquids", 2030, $null, $null, $null, 0.019593092998278099, 0.000081244902675383100, 0.0074870826105117002, 0.0092620082230198998, 0.00099351766891250001, 0.038036736199550301),
    @("Biomass [Solid]", 2030, $null, $null, 0.0033502898547493509, $null, $null, $null, $null, $null, $null),
    @("Renewable Energy Carrier", 2030, $null, $null, 0.00012473704652125531, $null, $null, $null, $null, $null, $null),
    @("Overall Demand", 2030, $null, $null, 0.0036042959845156270, 0.022140930841374279, 0.000092780936405462261, 0.0083062868836735956, 0.010788230023876770, 0.0011576540184250000, 0.041920026096714871),
    @("Hydrogen", 2040, $null, $null, $null, 0.00052863015891138902, $null, 0.000000025063597220873260, 0.000058505425606626571, $null, $null),
    @("Methanol", 2040, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @("Ammonia", 2040, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @("Synthetic Gases", 2040, $null, $null, $null, 0.00000000019296669613736470, $null, $null, 0.000000000019161140236206429, $null, $null),
    @("Biogenic Gases", 2040, $null, $null, 0.00052132502768251286, 0.000040820736754716933, $null, $null, 0.000013211566212401859, $null, $null),
    @("Fossil Gases", 2040, $null, $null, $null, 0.00022887893221654491, $null, $null, 0.000033244582423229073, $null, $null),
    @("Synthetic Liquids", 2040, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @("Biogenic Liquids", 2040, $null, $null, $null, 0.00085412606588417222, 0.000018805638135225190, 0.00099726376647110010, 0.00097286420463919994, 0.00019939902559690001, 0.0043936721089480454),
    @("Fossil Liquids", 2040, $null, $null, $null, 0.0052876814078324430, 0.000087342787361256887, 0.0070589198916128004, 0.0042151413393947002, 0.00088165804604870002, 0.036901038084350787),
    @("Biomass [Solid]", 2040, $null, $null, 0.0032744433964367798, $null, $null, $null, $null, $null, $null),
    @("Renewable Energy Carrier", 2040, $null, $null, 0.00050379097865586408, $null, $null, $null, $null, $null, $null),
    @("Overall Demand", 2040, $null, $null, 0.0042995594027751566, 0.0069401374945659618, 0.00010614842549648210, 0.0080562087216811216, 0.0052929671374372976, 0.0010810570716455999, 0.041294710193298827),
    @("Hydrogen", 2050, $null, $null, $null, 0.00073309596307136963, $null, 0.000000042480778456866487, 0.000093815099163969647, $null, $null),
    @("Methanol", 2050, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @("Ammonia", 2050, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @("Synthetic Gases", 2050, $null, $null, $null, 0.0000000017535708212727769, $null, $null, 0.00000000056674317233593355, $null, $null),
    @("Biogenic Gases", 2050, $null, $null, 0.0012874728085888040, 0.0000069222815955714644, $null, $null, 0.0000037003714757829182, $null, $null),
    @("Fossil Gases", 2050, $null, $null, $null, 0.000014009661013839189, $null, $null, 0.000011916743997400550, $null, $null),
    @("Synthetic Liquids", 2050, $null, $null, $null, 0.0000000000037928085978322656, 0.00000000000069251000880319691, 0.000000000038421863109496971, 0.000000000015590760160350768, 0.0000000000013656491720614300, 0.00000000030065018212667119),
    @("Biogenic Liquids", 2050, $null, $null, $null, 0.000078158561916318495, 0.000033602373614937248, 0.0013062006708533999, 0.00025213983144080002, 0.00025645872277050002, 0.0062593991433204359),
    @("Fossil Liquids", 2050, $null, $null, $null, 0.00025361655629910858, 0.000078779495201921284, 0.0064015315948015997, 0.00074273580111250002, 0.00075888795447479995, 0.034401809486705709),
    @("Biomass [Solid]", 2050, $null, $null, 0.0033867483194102281, $null, $null, $null, $null, $null, $null),
    @("Renewable Energy Carrier", 2050, $null, $null, 0.0012414394072420740, $null, $null, $null, $null, $null, $null),
    @("Overall Demand", 2050, $null, $null, 0.0059156605352411063, 0.0010858047812598371, 0.00011238186950936861, 0.0077077747848553202, 0.0011043084295243851, 0.0010153466786109490, 0.040661208930676332)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        if ($val -ne $null) {
            $ws.Cells.Item($r, $c).Value = $val
        }
        $c = $c + 1
    }
    $r = $r + 1
}
